$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 4 (9906c219-... file) handoff/handback datetimes
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-23 18:53:04"
$wsZhCn.Range("H4").Value = "2016-03-23 18:53:31"

# de-de sheet: row 4 (9906c219-... file) handoff/handback datetimes
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-23 18:53:10"
$wsDeDe.Range("H4").Value = "2016-03-23 18:53:37"
